# Slide 2 ("Important points") - readme/tutorial wording tweak:
#   * new paragraph "Unzip the Tools.zip in" added above the path line
#   * "If you don't follow my instruction, I will take points off."
#     paragraph removed
#   * "C:\Graphics\Tools" line split into two runs: "C:\Graphics" + "\Tools"

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tr  = $shp.TextFrame.TextRange

# Paragraphs (1-based) before editing:
#   4: Every students must installed all required libraries on a specific directory
#   5: C:\Graphics\Tools                              (lvl 1)
#   6: If you don't follow my instruction, I will take points off.   (lvl 2)
#   7: <empty>                                         (lvl 1)

# 1) Drop the "If you don't follow my instruction..." paragraph entirely.
$tr.Paragraphs(6,1).Delete()

# 2) Drop the now-trailing empty paragraph. Because it is the very last
#    paragraph in the text frame, removing it folds its endParaRPr onto the
#    preceding "C:\Graphics\Tools" paragraph instead of just vanishing.
$tr.Paragraphs(6,1).Delete()

# 3) Split "C:\Graphics\Tools" into two runs: "C:\Graphics" and "\Tools".
$dirPara   = $tr.Paragraphs(5,1)
$firstPart = $tr.Characters($dirPara.Start, 11)
$firstPart.Text = "C:\Graphics"

# 4) Insert the new "Unzip the Tools.zip in" paragraph right after the
#    "Every students..." paragraph; it picks up that paragraph's (level 0)
#    formatting, matching the target (no explicit <a:pPr>).
$tr.Paragraphs(4,1).InsertAfter("`rUnzip the Tools.zip in") | Out-Null
